{"js": "const body = context.document.body;\n\n// Locate the paragraph containing the website URL (robust to ordering).\nconst paras = body.paragraphs;\nparas.load(\"items,text\");\nawait context.sync();\n\nlet websitePara = null;\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text.indexOf('\"website\"') !== -1) {\n    websitePara = paras.items[i];\n    break;\n  }\n}\n\nif (!websitePara) {\n  throw new Error(\"website paragraph not found\");\n}\n\n// 1. Replace the old URL text with the new one inside that paragraph.\nconst results = websitePara.search(\"https://tokenpocket.pro\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"http://www.dlsj.xyy/\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2. Move the _GoBack bookmark to the end of the (now edited) website paragraph,\n//    mirroring Word's behaviour of relocating _GoBack to the last edit site.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nwebsitePara.getRange(\"End\").insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Replace the old website URL with the new one.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"https://tokenpocket.pro\"\n$find.Replacement.Text = \"http://www.dlsj.xyy/\"\n$find.Forward = $true\n$find.Wrap = 0\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# 2. Move the _GoBack bookmark to the end of the website paragraph's text (Word relocates\n#    _GoBack to the site of the most recent edit). Locate that paragraph by content so this\n#    keeps working if the document is restructured.\n$websitePara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like '*\"website\"*') {\n        $websitePara = $p\n        break\n    }\n}\n\nif ($websitePara -ne $null) {\n    if ($d.Bookmarks.Exists(\"_GoBack\")) {\n        $d.Bookmarks.Item(\"_GoBack\").Delete()\n    }\n\n    # Find the character offset right before the paragraph mark (end of the real text).\n    $probe = $websitePara.Range\n    $probe.MoveEnd(1, -1) | Out-Null\n    $textEnd = $probe.End\n\n    # Work around a collapsed-range quirk that misfires exactly at \"paragraph end minus one\"\n    # (the offset we need): temporarily append a one-character placeholder so the target\n    # offset is no longer a boundary case, add the bookmark there, then remove the placeholder.\n    $appendPoint = $d.Range($textEnd, $textEnd)\n    $appendPoint.InsertAfter(\"X\")\n\n    $target = $d.Range($textEnd, $textEnd)\n    $d.Bookmarks.Add(\"_GoBack\", $target) | Out-Null\n\n    $placeholder = $d.Range($textEnd, $textEnd + 1)\n    $placeholder.Delete()\n}\n"}
